# Updates cryptocurrency price/volume data in the "cryptos" worksheet.
# Generated from the authoritative cell-level diff; applies each changed
# cell value while forcing text semantics (matching the original inline
# string cells) so numeric-looking values like "116.22" or "51.566.09"
# are not silently reinterpreted as numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "51.566.09"
Set-TextCell $ws.Range("E2") "  +3.37%  "
Set-TextCell $ws.Range("D3") "2.746.46"
Set-TextCell $ws.Range("E3") "  +2.76%  "
Set-TextCell $ws.Range("E4") "  -0.02%  "
Set-TextCell $ws.Range("D5") "116.22"
Set-TextCell $ws.Range("E5") "  +2.43%  "
Set-TextCell $ws.Range("D6") "333.45"
Set-TextCell $ws.Range("E6") "  +2.33%  "
Set-TextCell $ws.Range("D7") "0.532"
Set-TextCell $ws.Range("E7") "  +0.83%  "
Set-TextCell $ws.Range("E8") "  +0.03%  "
Set-TextCell $ws.Range("D9") "0.577"
Set-TextCell $ws.Range("E9") "  +4.44%  "
Set-TextCell $ws.Range("E10") "  +1.81%  "
Set-TextCell $ws.Range("E11") "  +0.40%  "
Set-TextCell $ws.Range("D12") "0.0829"
Set-TextCell $ws.Range("E12") "  +0.93%  "
Set-TextCell $ws.Range("D14") "7.62"
Set-TextCell $ws.Range("E14") "  +3.54%  "
Set-TextCell $ws.Range("D15") "3.173.55"
Set-TextCell $ws.Range("E15") "  +2.61%  "
Set-TextCell $ws.Range("D16") "2.744.32"
Set-TextCell $ws.Range("E16") "  +2.43%  "
Set-TextCell $ws.Range("E17") "  +1.68%  "
Set-TextCell $ws.Range("D18") "51.504.70"
Set-TextCell $ws.Range("E18") "  +3.31%  "
Set-TextCell $ws.Range("D19") "13.78"
Set-TextCell $ws.Range("E19") "  +4.90%  "
Set-TextCell $ws.Range("E20") "  +2.17%  "
Set-TextCell $ws.Range("D21") "6.85"
Set-TextCell $ws.Range("E21") "  +1.02%  "
Set-TextCell $ws.Range("E22") "  +0.13%  "
Set-TextCell $ws.Range("D23") "276.92"
Set-TextCell $ws.Range("E23") "  +0.30%  "
Set-TextCell $ws.Range("D24") "70.23"
Set-TextCell $ws.Range("E24") "  -2.19%  "
Set-TextCell $ws.Range("D25") "2.68"
Set-TextCell $ws.Range("E25") "  +4.59%  "
Set-TextCell $ws.Range("D26") "26.90"
Set-TextCell $ws.Range("E26") "  +0.31%  "
Set-TextCell $ws.Range("E27") "  +0.50%  "
Set-TextCell $ws.Range("E28") "  +0.13%  "
Set-TextCell $ws.Range("D29") "10.30"
Set-TextCell $ws.Range("E29") "  +0.78%  "
Set-TextCell $ws.Range("E30") "  -1.11%  "
Set-TextCell $ws.Range("D31") "35.54"
Set-TextCell $ws.Range("E31") "  -1.57%  "
Set-TextCell $ws.Range("E32") "  +0.26%  "
Set-TextCell $ws.Range("D33") "50.33"
Set-TextCell $ws.Range("E33") "  +0.14%  "
Set-TextCell $ws.Range("D34") "5.62"
Set-TextCell $ws.Range("E34") "  +2.46%  "
Set-TextCell $ws.Range("D35") "0.0823"
Set-TextCell $ws.Range("E35") "  +2.16%  "
Set-TextCell $ws.Range("D36") "19.37"
Set-TextCell $ws.Range("E36") "  -0.66%  "
Set-TextCell $ws.Range("E37") "  -0.21%  "
Set-TextCell $ws.Range("D38") "2.11"
Set-TextCell $ws.Range("E38") "  +1.89%  "
Set-TextCell $ws.Range("E39") "  +5.10%  "
Set-TextCell $ws.Range("D40") "4.99"
Set-TextCell $ws.Range("E40") "  -0.81%  "
Set-TextCell $ws.Range("D41") "129.51"
Set-TextCell $ws.Range("E41") "  +3.13%  "
Set-TextCell $ws.Range("D42") "23.63"
Set-TextCell $ws.Range("E42") "  +5.95%  "
Set-TextCell $ws.Range("D43") "0.0348"
Set-TextCell $ws.Range("E43") "  +9.83%  "
Set-TextCell $ws.Range("B44") "WEMIXToken"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws.Range("D44") "2.29"
Set-TextCell $ws.Range("E44") "  +3.55%  "
Set-TextCell $ws.Range("B45") "Stellar"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D45") "0.113"
Set-TextCell $ws.Range("E45") "  +0.73%  "
Set-TextCell $ws.Range("D46") "2.37"
Set-TextCell $ws.Range("E46") "  +13.36%  "
Set-TextCell $ws.Range("D47") "2.101.47"
Set-TextCell $ws.Range("D48") "3.38"
Set-TextCell $ws.Range("E48") "  +1.98%  "
Set-TextCell $ws.Range("D49") "2.26"
Set-TextCell $ws.Range("E49") "  +1.91%  "
Set-TextCell $ws.Range("E50") "  +5.40%  "
Set-TextCell $ws.Range("D51") "8.97"
Set-TextCell $ws.Range("E51") "  -0.69%  "
